$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Safe (non-numeric-looking) text updates: set directly ---
# (includes the row 46/47 Coin/Link swap: ImmutableX <-> Stacks)
$ws.Range("D2").Value = "76.383.56"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.968.68"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  +5.34%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").Value = "2.966.05"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "3.515.35"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("E15").Value = "  +5.95%  "
$ws.Range("D16").Value = "76.332.18"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "2.971.07"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "3.121.21"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +5.82%  "
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +10.95%  "
$ws.Range("E40").Value = "  +14.93%  "
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("E48").Value = "  +7.93%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("E51").Value = "  +2.50%  "

# --- Numeric-looking Price text: force text format, set value, then restore
#     default ("Normal") cell style so no stray formatting is introduced.
#     (Plain .Value assignment would let Excel re-interpret strings like
#     "199.90" / "1.00" / "0.0000106" as numbers and mangle their text.)
$numericTextCells = @(
    @{ Ref = "D5"; Val = "199.90" },
    @{ Ref = "D6"; Val = "629.93" },
    @{ Ref = "D8"; Val = "0.548" },
    @{ Ref = "D9"; Val = "0.200" },
    @{ Ref = "D13"; Val = "4.98" },
    @{ Ref = "D15"; Val = "28.97" },
    @{ Ref = "D19"; Val = "13.39" },
    @{ Ref = "D20"; Val = "8.74" },
    @{ Ref = "D21"; Val = "371.65" },
    @{ Ref = "D22"; Val = "2.26" },
    @{ Ref = "D23"; Val = "4.28" },
    @{ Ref = "D24"; Val = "72.65" },
    @{ Ref = "D26"; Val = "1.00" },
    @{ Ref = "D27"; Val = "4.28" },
    @{ Ref = "D28"; Val = "9.67" },
    @{ Ref = "D29"; Val = "0.0000106" },
    @{ Ref = "D33"; Val = "504.50" },
    @{ Ref = "D34"; Val = "1.93" },
    @{ Ref = "D36"; Val = "164.03" },
    @{ Ref = "D37"; Val = "20.24" },
    @{ Ref = "D39"; Val = "0.381" },
    @{ Ref = "D40"; Val = "0.104" },
    @{ Ref = "D41"; Val = "183.33" },
    @{ Ref = "D44"; Val = "42.91" },
    @{ Ref = "D45"; Val = "4.91" },
    @{ Ref = "D46"; Val = "1.23" },
    @{ Ref = "D47"; Val = "1.63" },
    @{ Ref = "D48"; Val = "0.712" },
    @{ Ref = "D50"; Val = "2.30" },
    @{ Ref = "D51"; Val = "3.82" }
)

foreach ($item in $numericTextCells) {
    $rng = $ws.Range($item.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
    $rng.Style = "Normal"
}
